# "Doing Updates for Financials" - refresh TDLAF yearly financial figures
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TDLAF")

# Income Statement -------------------------------------------------------
# Total Revenue (row 8)
$ws.Range("D8").Value = 2200
$ws.Range("E8").Value = 2900
$ws.Range("F8").Value = 1000
$ws.Range("G8").Value = 2300
$ws.Range("H8").Value = 600
$ws.Range("I8").Value = 4100
$ws.Range("J8").Value = 2500
$ws.Range("D9").Value = 300
$ws.Range("E9").Value = 700
$ws.Range("F9").Value = 300
$ws.Range("G9").Value = 700
$ws.Range("H9").Value = 200
$ws.Range("I9").Value = 1800
$ws.Range("J9").Value = 1600

# Gross Profit (row 10)
$ws.Range("D10").Value = 2000
$ws.Range("E10").Value = 2200
$ws.Range("F10").Value = 700
$ws.Range("G10").Value = 1600
$ws.Range("H10").Value = 400
$ws.Range("I10").Value = 2200

# Research Development (row 12)
$ws.Range("D12").Value = 1500
$ws.Range("E12").Value = 2200
$ws.Range("F12").Value = 1000
$ws.Range("G12").Value = 2000
$ws.Range("H12").Value = 600
$ws.Range("I12").Value = 1000
$ws.Range("J12").Value = 700

# Non Recurring (row 14) - switched from "NA" text to numeric 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0

# Total Operating Expenses (row 17)
$ws.Range("D17").Value = "NA"
$ws.Range("E17").Value = 5600
$ws.Range("F17").Value = 3000
$ws.Range("G17").Value = 4800
$ws.Range("H17").Value = 1900
$ws.Range("I17").Value = 5300
$ws.Range("J17").Value = 4000

# Operating Income or Loss (row 18)
$ws.Range("D18").Value = "NA"
$ws.Range("E18").Value = -2700
$ws.Range("F18").Value = -2000
$ws.Range("G18").Value = -2400
$ws.Range("H18").Value = -1300
$ws.Range("I18").Value = -1300
$ws.Range("J18").Value = -1500

# Total Other Income/Expenses Net (row 20)
$ws.Range("D20").Value = "NA"

# Earnings Before Interest And Taxes (row 21)
$ws.Range("I21").Value = -1200
$ws.Range("J21").Value = "NA"

# Income Before Tax (row 23)
$ws.Range("D23").Value = "NA"
$ws.Range("E23").Value = -2700
$ws.Range("F23").Value = -2000
$ws.Range("G23").Value = -2400
$ws.Range("H23").Value = -1300
$ws.Range("I23").Value = -1300
$ws.Range("J23").Value = -1500

# Income After Tax (row 26)
$ws.Range("D26").Value = -1800
$ws.Range("E26").Value = -2700
$ws.Range("F26").Value = -2000
$ws.Range("G26").Value = -2400
$ws.Range("H26").Value = -1300
$ws.Range("I26").Value = -1300
$ws.Range("J26").Value = -1500

# Net Income From Continuing Ops (row 27)
$ws.Range("D27").Value = -1500
$ws.Range("E27").Value = -2700
$ws.Range("F27").Value = -2000
$ws.Range("G27").Value = -2400
$ws.Range("H27").Value = -1300
$ws.Range("I27").Value = -1300
$ws.Range("J27").Value = -1500

# Discontinued Operations (row 29)
$ws.Range("D29").Value = -2500
$ws.Range("E29").Value = -1900
$ws.Range("F29").Value = 100
$ws.Range("G29").Value = -6300
$ws.Range("H29").Value = 0

# Other Items (row 32) - switched from numeric 0 to "NA" text
$ws.Range("D32").Value = "NA"

# Net Income (row 33)
$ws.Range("D33").Value = -4000
$ws.Range("E33").Value = -4600
$ws.Range("F33").Value = -1900
$ws.Range("G33").Value = -8800
$ws.Range("H33").Value = -1300
$ws.Range("I33").Value = -1300
$ws.Range("J33").Value = -1500

# Net Income Applicable To Common Shares (row 35)
$ws.Range("D35").Value = -4000
$ws.Range("E35").Value = -4600
$ws.Range("F35").Value = -1900
$ws.Range("G35").Value = -8800
$ws.Range("H35").Value = -1300
$ws.Range("I35").Value = -1300
$ws.Range("J35").Value = -1500

# Balance Sheet -----------------------------------------------------------
# Cash And Cash Equivalents (row 41)
$ws.Range("D41").Value = 5600
$ws.Range("E41").Value = 7500
$ws.Range("F41").Value = 9500
$ws.Range("G41").Value = 800

# Short Term Investments (row 42)
$ws.Range("I42").Value = 5700

# Net Receivables (row 43)
$ws.Range("H43").Value = 2200

# Total Current Assets (row 46)
$ws.Range("D46").Value = 8300
$ws.Range("E46").Value = 9200
$ws.Range("F46").Value = 13200
$ws.Range("G46").Value = 3200
$ws.Range("H46").Value = 6400
$ws.Range("I46").Value = 7300
$ws.Range("J46").Value = 4900

# Property Plant and Equipment (row 48)
$ws.Range("E48").Value = 2300

# Other Assets (row 52)
$ws.Range("D52").Value = 600
$ws.Range("E52").Value = 3200

# Total Assets (row 54)
$ws.Range("D54").Value = 11600
$ws.Range("E54").Value = 15700
$ws.Range("F54").Value = 17500
$ws.Range("G54").Value = 3600
$ws.Range("H54").Value = 6800
$ws.Range("I54").Value = 7500
$ws.Range("J54").Value = 5200

# Other Current Liabilities (row 59)
$ws.Range("I59").Value = 300

# Total Current Liabilities (row 60)
$ws.Range("E60").Value = 800
$ws.Range("H60").Value = 1100

# Total Liabilities (row 66)
$ws.Range("E66").Value = 800
$ws.Range("G66").Value = 1100

# Retained Earnings (row 72)
$ws.Range("D72").Value = -15300
$ws.Range("E72").Value = -10800
$ws.Range("F72").Value = -9200
$ws.Range("G72").Value = -109600
$ws.Range("H72").Value = -106400
$ws.Range("I72").Value = -105200
$ws.Range("J72").Value = -103700

# Total Stockholder Equity (row 76)
$ws.Range("D76").Value = 10300
$ws.Range("E76").Value = 14800
$ws.Range("F76").Value = 16400
$ws.Range("G76").Value = 2400
$ws.Range("H76").Value = 5700
$ws.Range("I76").Value = 6900
$ws.Range("J76").Value = 4000

# Cash Flow Statement ------------------------------------------------------
# Net Income (row 81)
$ws.Range("D81").Value = -4000
$ws.Range("E81").Value = -4600
$ws.Range("F81").Value = -1900
$ws.Range("G81").Value = -8800
$ws.Range("H81").Value = -1300
$ws.Range("I81").Value = -1300
$ws.Range("J81").Value = -1500

# Depreciation (row 83) - switched from numeric 100 to "NA" text
$ws.Range("J83").Value = "NA"

# Changes In Inventories (row 89)
$ws.Range("D89").Value = -1700
$ws.Range("F89").Value = -3100
$ws.Range("G89").Value = -1200
$ws.Range("I89").Value = -1800
$ws.Range("J89").Value = -3800

# Total Cash Flow From Operating Activities (row 91)
$ws.Range("G91").Value = -2700

# Other Cashflows from Investing Activities (row 94)
$ws.Range("E94").Value = -1000
$ws.Range("F94").Value = 8600
$ws.Range("I94").Value = -2800
$ws.Range("J94").Value = "NA"

# Other Cash Flows from Financing Activities (row 100)
$ws.Range("G100").Value = -1900
$ws.Range("I100").Value = 4400
$ws.Range("J100").Value = "NA"

# Effect Of Exchange Rate Changes (row 101)
$ws.Range("G101").Value = 1700
$ws.Range("J101").Value = "NA"

# Change In Cash and Cash Equivalents (row 102)
$ws.Range("D102").Value = -1900
$ws.Range("F102").Value = 5700
$ws.Range("G102").Value = -3000
$ws.Range("J102").Value = -3200
